$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Insert "storeKeys(json,jsonpath,var)" into the json function list (column M),
#    pushing storeValue/storeValues down by one row (M16 -> M17 -> M18).
$m16 = $ws.Range("M16").Value()
$m17 = $ws.Range("M17").Value()
$ws.Range("M18").Value = $m17
$ws.Range("M17").Value = $m16
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# 2) Remove "text" from the master list of names (column A, row 25), shifting
#    webalert/webcookie/ws/ws.async/xml up by one row, and clearing the last row.
for ($r = 25; $r -le 30; $r++) {
    $nextVal = $ws.Range("A" + ($r + 1)).Value()
    $ws.Range("A" + $r).Value = $nextVal
}
$ws.Range("A31").ClearContents()

# 3) Delete column Y (the old, now-orphaned "text" column) so that web / webalert /
#    webcookie / ws / ws.async / xml shift one column to the left (Z->Y, AA->Z, ...).
$ws.Columns("Y").Delete()

# 4) Update the defined names to reflect the new ranges.
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"

Write-Output "done"
